$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the worksheet tab: "Hoja1" -> "Semilla 4"
# ---------------------------------------------------------------------------
$ws.Name = "Semilla 4"

# ---------------------------------------------------------------------------
# 2. Propagate existing cell formats to the new cells/columns that need them
#    BEFORE we overwrite the source cells' own values, so every destination
#    reuses an existing style (no spurious/duplicate style entries).
#    Style map (based on the pre-edit workbook):
#      style "2" source -> A1   (header style, centered)
#      style "3" source -> A2   (hyperlink-ish row style)
#      style "4" source -> D2   (hyperlink-ish row style, vertical center)
#      style "5" source -> A4   (data row style, center+vcenter)
#      style "1" source -> C2   (plain text style)
#      style "6" (new)  -> built from A8's old style (vertical-center), then
#                          re-colored/bolded into the new bold-blue font.
# ---------------------------------------------------------------------------

# style 2 -> row 3 and row 8 (A:C)
$ws.Range("A1").Copy()
$ws.Range("A3:D3").PasteSpecial(-4122)
$ws.Range("A8:C8").PasteSpecial(-4122)

# style 3 -> F... (not needed further, already used on A2:C2,E2)

# style 4 -> D8 (reuse D2's style)
$ws.Range("D2").Copy()
$ws.Range("D8").PasteSpecial(-4122)

# style 5 -> rows 5,6,7 already have it; nothing more needed

# style 1 -> F1, G1, and rows 9-13 (A:D)
$ws.Range("C2").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$ws.Range("A9:D13").PasteSpecial(-4122)

# Build the new bold/blue font style (style "6") in a scratch cell from the
# old A8 style (vertical-center, JetBrains Mono), then tweak its font.
$ws.Range("A8").Copy()
$ws.Range("Z100").PasteSpecial(-4122)
$ws.Range("Z100").Font.Bold = $true
$ws.Range("Z100").Font.Color = 14580521
$ws.Range("Z100").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("F2:H2").PasteSpecial(-4122)
$ws.Range("Z100").Clear()

# ---------------------------------------------------------------------------
# 3. Write the new cell values (text, so they keep the "@" text format
#    already present on every style used above).
# ---------------------------------------------------------------------------

# Row 1 (headers)
$ws.Range("A1").Value = "URL EPOS"
$ws.Range("B1").Value = "URL CRM"
$ws.Range("C1").Value = "URL CONFIRMADOR"
$ws.Range("D1").Value = "URL GATEWAYCBS"
$ws.Range("E1").Value = "URL GATEWAY MG"
$ws.Range("F1").Value = "user"
$ws.Range("G1").Value = "password"
$ws.Range("H1").Value = "rutaWinWap"

# Row 2 (values/URLs + new WinWap columns)
$ws.Range("A2").Value = "http://10.69.60.77:8180/tigo-pos-web/"
$ws.Range("B2").Value = "http://10.69.60.85:8280/portal/login?initialURI=%2Fportal%2FCRMPortal"
$ws.Range("C2").Value = "http://10.69.60.77:8180/tigo-pos-web/wap/windex.wml"
$ws.Range("D2").Value = "http://10.65.45.12:9001/gatewaycbs/BcServicesInt"
$ws.Range("E2").Value = "http://10.65.45.12:9001/gatewaymgint/GatewayMGWSInt"
$ws.Range("F2").Value = "CQ10960370"
$ws.Range("G2").Value = "Tigo.2022*"
$ws.Range("H2").Value = 'C:\Program Files (x86)\Winwap Technologies\WinWAP for Windows 4.2\WinWAP4.exe'

# Row 3 (DB section header)
$ws.Range("A3").Value = "URL DB"
$ws.Range("B3").Value = "service"
$ws.Range("C3").Value = "user"
$ws.Range("D3").Value = "password"

# Row 4-7 (DB rows, unchanged content)
$ws.Range("A4").Value = "10.69.60.89"
$ws.Range("B4").Value = "DEV11G"
$ws.Range("C4").Value = "ACTIVATOR"
$ws.Range("D4").Value = "ACTIVATOR"

$ws.Range("A5").Value = "10.69.60.88"
$ws.Range("B5").Value = "desepos"
$ws.Range("C5").Value = "epos"
$ws.Range("D5").Value = "epos"

$ws.Range("A6").Value = "10.69.60.88"
$ws.Range("B6").Value = "dev10g"
$ws.Range("C6").Value = "postsale"
$ws.Range("D6").Value = "postsale"

$ws.Range("A7").Value = "10.65.32.74"
$ws.Range("B7").Value = "siebel04"
$ws.Range("C7").Value = "siebel"
$ws.Range("D7").Value = "siebel"

# Row 8 (new "Vendedor / Cedula Cliente / MSIDN / MSI" header block)
$ws.Range("A8").Value = "Vendedor"
$ws.Range("B8").Value = "Cedula Cliente"
$ws.Range("C8").Value = "MSIDN"
$ws.Range("D8").Value = "MSI"

# Rows 9-13 (seed data)
$ws.Range("A9").Value = "10960370"
$ws.Range("B9").Value = "984108505"
$ws.Range("C9").Value = "3016875982"
$ws.Range("D9").Value = "732111198172291"

$ws.Range("A10").Value = "10960370"
$ws.Range("B10").Value = "835244140"
$ws.Range("C10").Value = "3016877591"
$ws.Range("D10").Value = "732111198172292"

$ws.Range("A11").Value = "10960370"
$ws.Range("B11").Value = "667299000"
$ws.Range("C11").Value = "3016875982"
$ws.Range("D11").Value = "732111198172291"

$ws.Range("A12").Value = "10960370"
$ws.Range("B12").Value = "835244140"
$ws.Range("C12").Value = "3016877411"
$ws.Range("D12").Value = "732111198172294"

$ws.Range("A13").Value = "10960370"
$ws.Range("B13").Value = "311615530"
$ws.Range("C13").Value = "3016876876"
$ws.Range("D13").Value = "732111198172293"

# ---------------------------------------------------------------------------
# 4. Add the new hyperlink on C2 (mirrors the ones already on A2/B2/D2/E2).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("C2"), "http://10.69.60.77:8180/tigo-pos-web/wap/windex.wml")

# ---------------------------------------------------------------------------
# 5. Match the saved selection/active cell.
# ---------------------------------------------------------------------------
$ws.Range("C13").Select()
